$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.665.97'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '3.797.46'
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").Value = '4.442.37'
$ws.Range("E14").Value = '  +1.00%  '
$ws.Range("D15").Value = '3.812.11'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.74%  '
$ws.Range("D17").Value = '67.696.02'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.93%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.83%  '
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.38%  '
$ws.Range("E26").Value = '  -0.49%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").Value = '3.943.38'
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.94%  '
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.06'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("D36").Value = '3.740.60'
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("E37").Value = '  +0.60%  '
$ws.Range("E38").Value = '  +2.42%  '
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.88%  '
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("E48").Value = '  +8.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '148.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '396.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("E51").Value = '  +11.12%  '
